$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "SmallTrackDiameter" row (row 4)
$ws.Rows(4).Delete()

# Remove the "LargeTrackRatio" row (now shifted up to row 4)
$ws.Rows(4).Delete()

# Remove the "ShallowTrackDepth" row (this was row 9 originally,
# now shifted up to row 7 after the two prior deletions). This row
# held the formula =B1+(2*B3).
$ws.Rows(7).Delete()

# Update OpticalFiberDiameter value (row 1)
$ws.Range("B1").Value2 = 1

# Update ScintillatorWidth value (now row 5) 145 -> 200
$ws.Range("B5").Value2 = 200

# Update ScintillatorHeight (now row 6): replace formula with plain
# value 200 (no 3d cut formula anymore)
$ws.Range("B6").Value2 = 200

# Update LoopLargeTrackRatio value (now row 8) 1 -> 1.5
$ws.Range("B8").Value2 = 1.5

# Leave the active selection on B10 (FiberSensorClerence value), matching
# where the author was last working in the sheet
[void]$ws.Range("B10").Select()
